# Apply "elapsed time y cpu" changes:
# - Add two new header columns: G1 "Elapsed Time", H1 "CPU" (same style as other headers)
# - Update B2, C2, D2 with refreshed metric values
# - Update F2 pipeline description text (now includes n_estimators=150, wrapped across two lines)
# - Add new G2 (elapsed time) and H2 (cpu) values
# - Dimension grows from A1:F2 to A1:H2 (handled automatically by the engine)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: copy formatting from an existing header cell (A1) then set text,
# so they pick up the bold/centered/bordered header style (style index 1).
$ws.Range("A1").Copy($ws.Range("G1"))
$ws.Range("A1").Copy($ws.Range("H1"))
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Updated metric values in row 2
$ws.Range("B2").Value = 0.5379381221575794
$ws.Range("C2").Value = 0.9892885096110098
$ws.Range("D2").Value = 0.6104905280394304

# Updated model description (now wraps across two lines with n_estimators=150)
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=5, n_estimators=150))])"

# New Elapsed Time / CPU values
$ws.Range("G2").Value = 0.1228586025167412
$ws.Range("H2").Value = 0.991
